$d = $word.ActiveDocument

$replacements = @(
    @("2025-10-12 Sunday", "2025-10-13 Monday"),
    @("73×80=", "98×34="),
    @("69×61=", "59×65="),
    @("26×21=", "28×55="),
    @("21×42=", "26×78="),
    @("51×41=", "90×48="),
    @("85×66=", "21×54="),
    @("83×42=", "18×43="),
    @("31×61=", "72×34="),
    @("92×50=", "31×11="),
    @("86×53=", "40×57="),
    @("26×32=", "48×98="),
    @("46×87=", "70×92="),
    @("57×68=", "64×63="),
    @("63×20=", "21×92="),
    @("71×75=", "45×80="),
    @("11×36=", "67×59="),
    @("27×38=", "35×79="),
    @("93×70=", "71×90="),
    @("62×64=", "71×60="),
    @("23×74=", "92×54="),
    @("27×40=", "97×23="),
    @("85×24=", "79×98="),
    @("56×67=", "23×85="),
    @("20×31=", "36×53="),
    @("30×67=", "79×13=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
